# Updates cryptos list values (price/volume) and re-orders two
# coin rows (Monero/WstETH and Decentraland/EnergySwap) to match
# the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.803.49'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '1.848.06'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''335.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').Value = '''1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = '''0.4655'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('D8').Value = '''0.3868'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = '''46.75'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').Value = '''0.07915'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Value = '''0.9696'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = '''21.31'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').Value = '1.844.53'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = '''5.904'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').Value = '''7.157'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '''90.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '''0.06618'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = '''0.00001030'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').Value = '''17.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').Value = '''1.007'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '27.816.36'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').Value = '''5.352'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '''2.298'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').Value = '2.074.03'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''158.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').Value = '''19.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').Value = '''2.070'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('D30').Value = '''5.388'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').Value = '''118.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').Value = '''0.09427'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = '''0.9450'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').Value = '''3.596'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('D35').Value = '''5.263'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').Value = '''1.330'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.48%  '
$ws.Range('D37').Value = '''0.06024'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '''0.02217'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '''8.242'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.54%  '
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').Value = '''1.154'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').Value = '''0.5821'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '''0.1849'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').Value = '''1.281'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5459'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''11.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').Value = '''1.943'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').Value = '''0.06856'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('D50').Value = '''110.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
